$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A61").Value = "'2025-11-29"
$ws.Range("A61").Style = "Normal"
$ws.Range("B61").Value = 500
$ws.Range("C61").Value = 500
$ws.Range("D61").Value = 50
$ws.Range("E61").Value = 100
